$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "In Translation"
#    Appears in the "Overview" sheet (columns E & F, rows 2-3) and in each
#    per-locale sheet's "Status" column (column C, rows 2-3).
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2. Narrow the "zh-cn"/"de-de" status columns (Overview!E:F, zh-cn!C,
#    de-de!C) from their old width down to match the new, shorter status
#    text ("In Translation" vs "Ready for handoff").
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
